$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Row 2
$ws.Range("G2").Value = 0.01880741119384766
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.01880741119384766

# Row 3
$ws.Range("G3").Value = 0.01897311210632324
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.01897311210632324

# Row 4
$ws.Range("G4").Value = 0.02179503440856934
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.02179503440856934
